# Rename the app from "Grade.ly" to "Graded":
#  - app_name row (row 8, columns B:E): "Grade.ly" -> "Graded"
#  - github_summary row (row 43, columns B:E): repo URL updated to .../Graded
#    and turned into a hyperlink (as Excel does when you paste/confirm a URL)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("string")

# --- app_name row ---------------------------------------------------------
$ws.Range("B8:E8").Value = "Graded"

# --- github_summary row ----------------------------------------------------
$newUrl = "https://github.com/NightDreamGames/Graded"

# Update the text first so every cell already shows the new URL.
$ws.Range("B43:E43").Value = $newUrl

# Turn B43 into its own hyperlink (no custom display text).
$ws.Hyperlinks.Add($ws.Range("B43"), $newUrl)

# Turn C43:E43 into a single hyperlink spanning the three cells, with an
# explicit display string (mirrors how Excel records a dragged-fill of a
# hyperlinked cell across a range).
$ws.Hyperlinks.Add($ws.Range("C43:E43"), $newUrl, "", "", $newUrl)

# Make sure every cell in the row actually carries the Hyperlink cell style.
$ws.Range("B43:E43").Style = "Hyperlink"

# --- restore view state (best effort) --------------------------------------
$ws.Range("B43:E43").Select()
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
